$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8165242075920105
$ws.Range("B1").Value = 0.7612819671630859
$ws.Range("D1").Value = 1.527166604995728
$ws.Range("E1").Value = 0.9299272298812866
